# Card21: add a new service event (row 19) and backfill the "nan"
# placeholders on the previous event row (row 18) that were left as
# blank cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# --- Row 18: the columns that don't apply to this event were exported
#     as literal "nan" text (pandas/openpyxl style) everywhere else in
#     this table, except this last row where they were left blank.
#     Backfill them so row 18 matches the rest of the table.
$ws.Range("B18:K18").Value = "nan"
$ws.Range("M18").Value = "nan"

# --- Row 19: the new service event for Card21.
# Column A ("card") is "21" for every row in this table and is stored
# as text, not a number. Copy it down from A18 (instead of typing a
# fresh "21") so Excel keeps it as text instead of auto-converting the
# digits to a numeric value.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("L19").Value = "13\5\2025"
$ws.Range("M19").Value = "610.2 t"
$ws.Range("N19").Value = "تم عمل صيانه وتغيير الجرائد الاماميه (1_2_4_5_7_8) ومعايره المكنه"
$ws.Range("O19").Value = "الخبير"
